$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Producer Consumer row (row 6): code edit resolved the queue-overflow
# deadlock issue -> status flips from "needs fixing" to "OK", with the
# states/transitions reached now recorded.
$ws.Range("D6").Value = "OK"
$ws.Range("E6").Value = "msgsrv -> constructor, added else in giveNextProduce()"
$ws.Range("F6").Value = 161
$ws.Range("G6").Value = 275

# Reflect the improved status with the "Good" (green) conditional style,
# matching the other "OK" rows.
$ws.Range("D6:G6").Style = "Good"
